# Regenerate merged AHB files:
#  - rename the "_old" / "_new" header-column suffixes to "_FV2410" / "_FV2504"
#  - turn the data range into a native Excel Table ("Table1")
#  - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the 20 header cells (columns A-J = "_old" -> "_FV2410",
#        columns L-U = "_new" -> "_FV2504"; column K ("diff") is untouched) ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2. Turn A1:U89 into a real Table ("Table1") ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U89"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
